$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 4 and row 5 swap (file 62cc9c08 now reports a failed handback, file bf715205 moves down) ---
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(4,1).Value = "62cc9c08-8915-4986-abb9-e5c2a3f3b513.md"
$wsOverview.Cells.Item(4,2).Value = "Handback transform failed"
$wsOverview.Cells.Item(4,3).Value = "Handback transform failed"
$wsOverview.Cells.Item(4,4).Value = "2016-03-23 07:33:28"

$wsOverview.Cells.Item(5,1).Value = "bf715205-6621-48a5-afcd-4af86d10eaaf.md"
$wsOverview.Cells.Item(5,2).Value = "In Translation"
$wsOverview.Cells.Item(5,3).Value = "In Translation"
$wsOverview.Cells.Item(5,4).Value = "2016-03-23 07:30:57"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Cells.Item(4,1).Value = "62cc9c08-8915-4986-abb9-e5c2a3f3b513.md"
$wsZh.Cells.Item(4,3).Value = "Handback transform failed"
$wsZh.Cells.Item(4,4).Value = "62cc9c08-8915-4986-abb9-e5c2a3f3b513.980219e2f2a1c4736827b8f232569681a12960f8.zh-cn.xlf"
$wsZh.Cells.Item(4,5).Value = "2016-03-23 07:33:20"
$wsZh.Cells.Item(4,12).Value = "The handback type mt is not match with handoff type ht."

$wsZh.Cells.Item(5,1).Value = "bf715205-6621-48a5-afcd-4af86d10eaaf.md"
$wsZh.Cells.Item(5,3).Value = "In Translation"
$wsZh.Cells.Item(5,4).Value = "bf715205-6621-48a5-afcd-4af86d10eaaf.288b120a6ada4a457c530043517ea3473e231327.zh-cn.xlf"
$wsZh.Cells.Item(5,5).Value = "2016-03-23 07:30:48"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Cells.Item(4,1).Value = "62cc9c08-8915-4986-abb9-e5c2a3f3b513.md"
$wsDe.Cells.Item(4,3).Value = "Handback transform failed"
$wsDe.Cells.Item(4,4).Value = "62cc9c08-8915-4986-abb9-e5c2a3f3b513.980219e2f2a1c4736827b8f232569681a12960f8.de-de.xlf"
$wsDe.Cells.Item(4,5).Value = "2016-03-23 07:33:28"
$wsDe.Cells.Item(4,12).Value = "The handback type mt is not match with handoff type ht."

$wsDe.Cells.Item(5,1).Value = "bf715205-6621-48a5-afcd-4af86d10eaaf.md"
$wsDe.Cells.Item(5,3).Value = "In Translation"
$wsDe.Cells.Item(5,4).Value = "bf715205-6621-48a5-afcd-4af86d10eaaf.288b120a6ada4a457c530043517ea3473e231327.de-de.xlf"
$wsDe.Cells.Item(5,5).Value = "2016-03-23 07:30:57"
